$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the corrected/normalized user record.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "michael"
$ws.Range("C2").Value = "100-534"
$ws.Range("D2").Value = "Graham Gibson"
$ws.Range("E2").Value = "Computer Science"
$ws.Range("F2").Value = "Queen's University"
$ws.Range("G2").Value = "CMC"
$ws.Range("H2").Value = "Academic Machine Dependent"

# Permissions code is a numeric-looking string; keep it stored as text.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "1111111"

# Remove the now-duplicate/stale user rows (3 and 4).
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()
